$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.04870513615401251
$ws.Range("D2").Value = 0.179318368419942
$ws.Range("E2").Value = 0.1637399210225325
$ws.Range("F2").Value = 1.496597299282612
$ws.Range("G2").Value = 0.002457148500969595
$ws.Range("J2").Value = 0.1897977685385754
$ws.Range("K2").Value = 1.445185142343121
$ws.Range("N2").Value = 1.194661019790658
$ws.Range("O2").Value = 3.642274077026741

$ws.Range("C3").Value = 0.04327683841606245
$ws.Range("D3").Value = 0.1748351904357719
$ws.Range("E3").Value = 0.1599260690104245
$ws.Range("F3").Value = 1.492582903046312
$ws.Range("G3").Value = 0.002460491392391409
$ws.Range("J3").Value = 0.1855785131881618
$ws.Range("K3").Value = 1.304952777451831
$ws.Range("N3").Value = 1.200710615396503
$ws.Range("O3").Value = 3.643885872975886

$ws.Range("C4").Value = 0.03995888928245961
$ws.Range("D4").Value = 0.1721467448533076
$ws.Range("E4").Value = 0.1576655537974432
$ws.Range("F4").Value = 1.491071690240759
$ws.Range("G4").Value = 0.002462653368624975
$ws.Range("J4").Value = 0.1831004043892364
$ws.Range("K4").Value = 1.218980614162859
$ws.Range("N4").Value = 1.204908205562653
$ws.Range("O4").Value = 3.64735405205542

$ws.Range("C5").Value = 0.03861055267469737
$ws.Range("D5").Value = 0.17106741717609
$ws.Range("E5").Value = 0.1567648159037276
$ws.Range("F5").Value = 1.490695343703891
$ws.Range("G5").Value = 0.002463561992923302
$ws.Range("J5").Value = 0.1821188209089541
$ws.Range("K5").Value = 1.183980919902069
$ws.Range("N5").Value = 1.206740423388617
$ws.Range("O5").Value = 3.649389453026487

$ws.Range("C6").Value = 0.03838688833036485
$ws.Range("D6").Value = 0.1708891786311284
$ws.Range("E6").Value = 0.1566164842341884
$ws.Range("F6").Value = 1.490647306756941
$ws.Range("G6").Value = 0.002463714538999455
$ws.Range("J6").Value = 0.181957536028321
$ws.Range("K6").Value = 1.178171391041332
$ws.Range("N6").Value = 1.207052015588161
$ws.Range("O6").Value = 3.649764968775941

$ws.Range("C7").Value = 0.03994068997779721
$ws.Range("D7").Value = 0.1721321228170183
$ws.Range("E7").Value = 0.1576533233406394
$ws.Range("F7").Value = 1.491065645457141
$ws.Range("G7").Value = 0.002462665510701647
$ws.Range("J7").Value = 0.1830870520036925
$ws.Range("K7").Value = 1.218508453448834
$ws.Range("N7").Value = 1.204932422629135
$ws.Range("O7").Value = 3.647378984912024

$ws.Range("C8").Value = 0.04683031671685001
$ws.Range("D8").Value = 0.1777592816769271
$ws.Range("E8").Value = 0.1624080614180414
$ws.Range("F8").Value = 1.495015010387107
$ws.Range("G8").Value = 0.002458278469449693
$ws.Range("J8").Value = 0.1883195929998038
$ws.Range("K8").Value = 1.396806791395818
$ws.Range("N8").Value = 1.196646756614882
$ws.Range("O8").Value = 3.642314697492452

$ws.Range("C9").Value = 0.06046225882805345
$ws.Range("D9").Value = 0.1893010195657325
$ws.Range("E9").Value = 0.1723761941183781
$ws.Range("F9").Value = 1.510343410554867
$ws.Range("G9").Value = 0.002450539864555716
$ws.Range("J9").Value = 0.1994756642249484
$ws.Range("K9").Value = 1.74743252759788
$ws.Range("N9").Value = 1.184224641912678
$ws.Range("O9").Value = 3.652107115553179

$ws.Range("C10").Value = 0.07055543272190334
$ws.Range("D10").Value = 0.1980869757766044
$ws.Range("E10").Value = 0.1800931631233453
$ws.Range("F10").Value = 1.526255843182099
$ws.Range("G10").Value = 0.002445375803482976
$ws.Range("J10").Value = 0.2082218760261014
$ws.Range("K10").Value = 2.005586776834832
$ws.Range("N10").Value = 1.177422002625036
$ws.Range("O10").Value = 3.671412697767209

$ws.Range("C11").Value = 0.07516483727422951
$ws.Range("D11").Value = 0.2021499145936332
$ws.Range("E11").Value = 0.1836894536793423
$ws.Range("F11").Value = 1.534510856840413
$ws.Range("G11").Value = 0.002443138614636537
$ws.Range("J11").Value = 0.2123211242833491
$ws.Range("K11").Value = 2.123138345310224
$ws.Range("N11").Value = 1.174830311703673
$ws.Range("O11").Value = 3.682844455713195

$ws.Range("C12").Value = 0.07691292988793919
$ws.Range("D12").Value = 0.2036978923513999
$ws.Range("E12").Value = 0.1850636105569663
$ws.Range("F12").Value = 1.53778340021384
$ws.Range("G12").Value = 0.002442307462363735
$ws.Range("J12").Value = 0.2138907935299414
$ws.Range("K12").Value = 2.167667412368644
$ws.Range("N12").Value = 1.173921089045521
$ws.Range("O12").Value = 3.687555809851801

$ws.Range("C13").Value = 0.07653633030284368
$ws.Range("D13").Value = 0.2033640895187858
$ws.Range("E13").Value = 0.1847671136587223
$ws.Range("F13").Value = 1.537072075311258
$ws.Range("G13").Value = 0.002442485754520851
$ws.Range("J13").Value = 0.2135519639072783
$ws.Range("K13").Value = 2.158076650051271
$ws.Range("N13").Value = 1.174113697321687
$ws.Range("O13").Value = 3.686524106053042

$ws.Range("C14").Value = 0.07530860141834239
$ws.Range("D14").Value = 0.2022770791191704
$ws.Range("E14").Value = 0.1838022595932998
$ws.Range("F14").Value = 1.534777151618712
$ws.Range("G14").Value = 0.002443069914476536
$ws.Range("J14").Value = 0.2124499135219651
$ws.Range("K14").Value = 2.126801491484684
$ws.Range("N14").Value = 1.174754063165636
$ws.Range("O14").Value = 3.68322438987957

$ws.Range("C15").Value = 0.07455692288152704
$ws.Range("D15").Value = 0.2016124795953118
$ws.Range("E15").Value = 0.1832128624584541
$ws.Range("F15").Value = 1.53339054137659
$ws.Range("G15").Value = 0.002443429813699444
$ws.Range("J15").Value = 0.2117771394248251
$ws.Range("K15").Value = 2.107646432268439
$ws.Range("N15").Value = 1.175155704249917
$ws.Range("O15").Value = 3.681253059065767

$ws.Range("C16").Value = 0.07025456120760509
$ws.Range("D16").Value = 0.1978227774235108
$ws.Range("E16").Value = 0.1798598620487226
$ws.Range("F16").Value = 1.525736844387524
$ws.Range("G16").Value = 0.002445524255232034
$ws.Range("J16").Value = 0.2079564095442379
$ws.Range("K16").Value = 1.997906681192376
$ws.Range("N16").Value = 1.177601484469392
$ws.Range("O16").Value = 3.670719040228306

$ws.Range("C17").Value = 0.06761982339615713
$ws.Range("D17").Value = 0.1955148045728521
$ws.Range("E17").Value = 0.1778248676842296
$ws.Range("F17").Value = 1.521302146859981
$ws.Range("G17").Value = 0.002446837749477178
$ws.Range("J17").Value = 0.2056434185740983
$ws.Range("K17").Value = 1.9306132805022
$ws.Range("N17").Value = 1.179230601369113
$ws.Range("O17").Value = 3.664936364282852

$ws.Range("C18").Value = 0.06610608161648202
$ws.Range("D18").Value = 0.1941935521841174
$ws.Range("E18").Value = 0.1766624700003518
$ws.Range("F18").Value = 1.518847067690714
$ws.Range("G18").Value = 0.002447603780085514
$ws.Range("J18").Value = 0.2043243920676474
$ws.Range("K18").Value = 1.891918914404073
$ws.Range("N18").Value = 1.180214964715461
$ws.Range("O18").Value = 3.66185962671193

$ws.Range("C19").Value = 0.06559384403955448
$ws.Range("D19").Value = 0.1937472718150985
$ws.Range("E19").Value = 0.1762702899641013
$ws.Range("F19").Value = 1.518032234560494
$ws.Range("G19").Value = 0.002447864958307844
$ws.Range("J19").Value = 0.2038797396832308
$ws.Range("K19").Value = 1.87881961516041
$ws.Range("N19").Value = 1.180556387335116
$ws.Range("O19").Value = 3.660860672212749

$ws.Range("C20").Value = 0.06790012065762596
$ws.Range("D20").Value = 0.1957598477823552
$ws.Range("E20").Value = 0.1780406604795814
$ws.Range("F20").Value = 1.521764326570931
$ws.Range("G20").Value = 0.002446696835154287
$ws.Range("J20").Value = 0.2058884660050921
$ws.Range("K20").Value = 1.937775650824733
$ws.Range("N20").Value = 1.179052280614442
$ws.Range("O20").Value = 3.665526127819049

$ws.Range("C21").Value = 0.07566914407074421
$ws.Range("D21").Value = 0.2025961051295013
$ws.Range("E21").Value = 0.1840853264149587
$ws.Range("F21").Value = 1.53544724587411
$ws.Range("G21").Value = 0.002442897898260874
$ws.Range("J21").Value = 0.2127731406376228
$ws.Range("K21").Value = 2.135987375453908
$ws.Range("N21").Value = 1.174564013750171
$ws.Range("O21").Value = 3.684183207271701

$ws.Range("C22").Value = 0.08076188229091485
$ws.Range("D22").Value = 0.2071189272864729
$ws.Range("E22").Value = 0.1881076651864007
$ws.Range("F22").Value = 1.545244208908599
$ws.Range("G22").Value = 0.002440508430271089
$ws.Range("J22").Value = 0.2173739537052626
$ws.Range("K22").Value = 2.265615850196582
$ws.Range("N22").Value = 1.172051427655035
$ws.Range("O22").Value = 3.698606231249073

$ws.Range("C23").Value = 0.07804238934011209
$ws.Range("D23").Value = 0.2047000130121432
$ws.Range("E23").Value = 0.1859543036475273
$ws.Range("F23").Value = 1.539937077661065
$ws.Range("G23").Value = 0.002441775217948511
$ws.Range("J23").Value = 0.2149091335379865
$ws.Range("K23").Value = 2.196423457229344
$ws.Range("N23").Value = 1.173353979545382
$ws.Range("O23").Value = 3.690703926790434

$ws.Range("C24").Value = 0.06777339514232494
$ws.Range("D24").Value = 0.1956490461991507
$ws.Range("E24").Value = 0.1779430770464501
$ws.Range("F24").Value = 1.521555080824044
$ws.Range("G24").Value = 0.002446760508549294
$ws.Range("J24").Value = 0.2057776466036358
$ws.Range("K24").Value = 1.93453756329103
$ws.Range("N24").Value = 1.179132750635944
$ws.Range("O24").Value = 3.665258723925149

$ws.Range("C25").Value = 0.0567610695408689
$ws.Range("D25").Value = 0.1861246953917686
$ws.Range("E25").Value = 0.1696105209711476
$ws.Range("F25").Value = 1.505381957283248
$ws.Range("G25").Value = 0.002452541391114899
$ws.Range("J25").Value = 0.1963613961133603
$ws.Range("K25").Value = 1.652478913770153
$ws.Range("N25").Value = 1.184224641912678
$ws.Range("O25").Value = 3.652107115553179
